$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any residual formatting on G61 (previously held a lone styled blank cell)
$ws.Range("G61").ClearFormats()

# Row 61: new Bus / Petrol ICE capital cost entry, assumed same as Diesel ICE
$ws.Range("B61").Value = "Assumed same as Diesel ICE"
$ws.Range("C61").Value = 2025
$ws.Range("D61").Value = "Bus"
$ws.Range("E61").Value = "Bus"
$ws.Range("F61").Value = "Petrol ICE"
$ws.Range("G61").Value = 400000
$ws.Range("K61").Value = 0.86
$ws.Range("K61").NumberFormat = "0.00"

# Row 62: new Bus / LPG capital cost entry, assumed same as Diesel ICE
$ws.Range("B62").Value = "Assumed same as Diesel ICE"
$ws.Range("C62").Value = 2025
$ws.Range("D62").Value = "Bus"
$ws.Range("E62").Value = "Bus"
$ws.Range("F62").Value = "LPG"
$ws.Range("G62").Value = 400000
$ws.Range("K62").Value = 0.86
$ws.Range("K62").NumberFormat = "0.00"

# Extend the autofilter to cover the new rows
$ws.AutoFilterMode = $false
$ws.Range("A1:K63").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "AG_costs!_FilterDatabase") {
        $n.RefersTo = "=AG_costs!`$A`$1:`$K`$63"
    }
}
